$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
